# Apply cached-value refresh from the scheduled Aegis Profits data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 661535.25
$ws.Range("I39").Value = 881970.9
$ws.Range("J39").Value = 228.33333
$ws.Range("K39").Value = 2645912.7
$ws.Range("L39").Value = 684.99999
$ws.Range("M39").Value = -2645616.7
$ws.Range("N39").Value = -1276.99999
$ws.Range("H42").Value = 1786088.6
$ws.Range("I42").Value = 2500107.2
$ws.Range("J42").Value = 1042
$ws.Range("K42").Value = 7500321.600000001
$ws.Range("L42").Value = 3126
$ws.Range("M42").Value = -7500091.600000001
$ws.Range("N42").Value = -3586
$ws.Range("H43").Value = 1369.0769
$ws.Range("I43").Value = 400
$ws.Range("J43").Value = 1449.8334
$ws.Range("K43").Value = 400
$ws.Range("L43").Value = 1449.8334
$ws.Range("M43").Value = -331
$ws.Range("N43").Value = -1587.8334
$ws.Range("H88").Value = 2572.9546
$ws.Range("I88").Value = 515.1429000000001
$ws.Range("J88").Value = 3533.2666
$ws.Range("K88").Value = 515.1429000000001
$ws.Range("L88").Value = 3533.2666
$ws.Range("M88").Value = -109.1429000000001
$ws.Range("N88").Value = -4345.2666
$ws.Range("H91").Value = 2572.9546
$ws.Range("I91").Value = 515.1429000000001
$ws.Range("J91").Value = 3533.2666
$ws.Range("K91").Value = 515.1429000000001
$ws.Range("L91").Value = 3533.2666
$ws.Range("M91").Value = 888.8570999999999
$ws.Range("N91").Value = -6341.2666
$ws.Range("H92").Value = 806.61536
$ws.Range("I92").Value = 888
$ws.Range("J92").Value = 623.5
$ws.Range("K92").Value = 888
$ws.Range("L92").Value = 623.5
$ws.Range("M92").Value = 360
$ws.Range("N92").Value = -3119.5
$ws.Range("H107").Value = 418.84
$ws.Range("J107").Value = 538.75
$ws.Range("L107").Value = 538.75
$ws.Range("N107").Value = -4378.75
$ws.Range("H129").Value = 636887.75
$ws.Range("J129").Value = 695526.5600000001
$ws.Range("L129").Value = 2086579.68
$ws.Range("N129").Value = -2096579.68
$ws.Range("H132").Value = 8626107
$ws.Range("I132").Value = 9621120
$ws.Range("J132").Value = 2668.6667
$ws.Range("K132").Value = 28863360
$ws.Range("L132").Value = 8006.000100000001
$ws.Range("M132").Value = -28860830
$ws.Range("N132").Value = -13066.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 15833.667
$ws.Range("J6").Value = 9000
$ws.Range("L6").Value = 9000
$ws.Range("N6").Value = -9346
$ws.Range("H61").Value = 1255.9032
$ws.Range("I61").Value = 1172.1072
$ws.Range("K61").Value = 1172.1072
$ws.Range("M61").Value = -960.1071999999999
$ws.Range("H74").Value = 3167.4736
$ws.Range("I74").Value = 1622.75
$ws.Range("J74").Value = 5815.5713
$ws.Range("K74").Value = 1622.75
$ws.Range("L74").Value = 5815.5713
$ws.Range("M74").Value = -748.75
$ws.Range("N74").Value = -7563.5713
$ws.Range("H77").Value = 3167.4736
$ws.Range("I77").Value = 1622.75
$ws.Range("J77").Value = 5815.5713
$ws.Range("K77").Value = 8113.75
$ws.Range("L77").Value = 29077.8565
$ws.Range("M77").Value = -3745.75
$ws.Range("N77").Value = -37813.85649999999
$ws.Range("H98").Value = 11205.5
$ws.Range("J98").Value = 11205.5
$ws.Range("L98").Value = 11205.5
$ws.Range("N98").Value = -17195.5
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("H122").Value = 2879.2727
$ws.Range("I122").Value = 2356
$ws.Range("J122").Value = 3795
$ws.Range("K122").Value = 7068
$ws.Range("L122").Value = 11385
$ws.Range("M122").Value = -4618
$ws.Range("N122").Value = -16285
$ws.Range("H132").Value = 4374.9414
$ws.Range("I132").Value = 4934
$ws.Range("K132").Value = 14802
$ws.Range("M132").Value = -12272
$ws.Range("H136").Value = 1255.9032
$ws.Range("I136").Value = 1172.1072
$ws.Range("K136").Value = 3516.3216
$ws.Range("M136").Value = -966.3215999999998
$ws.Range("N118").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 502777.5
$ws.Range("I7").Value = 502777.5
$ws.Range("K7").Value = 502777.5
$ws.Range("M7").Value = -502664.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23301.729
$ws.Range("I31").Value = 51565.1
$ws.Range("J31").Value = 3113.6072
$ws.Range("K31").Value = 51565.1
$ws.Range("L31").Value = 3113.6072
$ws.Range("M31").Value = -51270.1
$ws.Range("N31").Value = -3703.6072
$ws.Range("H34").Value = 23301.729
$ws.Range("I34").Value = 51565.1
$ws.Range("J34").Value = 3113.6072
$ws.Range("K34").Value = 51565.1
$ws.Range("L34").Value = 3113.6072
$ws.Range("M34").Value = -51363.1
$ws.Range("N34").Value = -3517.6072
$ws.Range("H45").Value = 14166.667
$ws.Range("H58").Value = 20337.572
$ws.Range("I58").Value = 2212.4443
$ws.Range("J58").Value = 52962.8
$ws.Range("K58").Value = 2212.4443
$ws.Range("L58").Value = 52962.8
$ws.Range("M58").Value = -2009.4443
$ws.Range("N58").Value = -53368.8
$ws.Range("H132").Value = 36588290
$ws.Range("I132").Value = 37039744
$ws.Range("K132").Value = 111119232
$ws.Range("M132").Value = -111116702
$ws.Range("H134").Value = 1621.4584
$ws.Range("I134").Value = 1757.25
$ws.Range("K134").Value = 5271.75
$ws.Range("M134").Value = -2736.75
$ws.Range("H136").Value = 20337.572
$ws.Range("I136").Value = 2212.4443
$ws.Range("J136").Value = 52962.8
$ws.Range("K136").Value = 6637.3329
$ws.Range("L136").Value = 158888.4
$ws.Range("M136").Value = -4087.3329
$ws.Range("N136").Value = -163988.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2334
$ws.Range("I139").Value = 1303.6364
$ws.Range("J139").Value = 3143.5715
$ws.Range("K139").Value = 3910.9092
$ws.Range("L139").Value = 9430.7145
$ws.Range("M139").Value = 1229.0908
$ws.Range("N139").Value = -19710.7145
$ws.Range("H141").Value = 2935.4375
$ws.Range("I141").Value = 2244
$ws.Range("J141").Value = 4456.6
$ws.Range("K141").Value = 6732
$ws.Range("L141").Value = 13369.8
$ws.Range("M141").Value = -1552
$ws.Range("N141").Value = -23729.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 8334483
$ws.Range("I3").Value = 10000960
$ws.Range("J3").Value = 2100
$ws.Range("K3").Value = 10000960
$ws.Range("L3").Value = 2100
$ws.Range("M3").Value = -10000844
$ws.Range("N3").Value = -2332
$ws.Range("H80").Value = 71430584
$ws.Range("I80").Value = 125001760
$ws.Range("J80").Value = 2346.6667
$ws.Range("K80").Value = 125001760
$ws.Range("L80").Value = 2346.6667
$ws.Range("M80").Value = -125000762
$ws.Range("N80").Value = -4342.6667
$ws.Range("H83").Value = 71430584
$ws.Range("I83").Value = 125001760
$ws.Range("J83").Value = 2346.6667
$ws.Range("K83").Value = 625008800
$ws.Range("L83").Value = 11733.3335
$ws.Range("M83").Value = -625003808
$ws.Range("N83").Value = -21717.3335
$ws.Range("H97").Value = 71430980
$ws.Range("I97").Value = 76925580
$ws.Range("J97").Value = 980
$ws.Range("K97").Value = 76925580
$ws.Range("L97").Value = 980
$ws.Range("M97").Value = -76925084
$ws.Range("N97").Value = -1972
$ws.Range("H102").Value = 263129.7
$ws.Range("I102").Value = 2240.2666
$ws.Range("J102").Value = 752297.4
$ws.Range("K102").Value = 2240.2666
$ws.Range("L102").Value = 752297.4
$ws.Range("M102").Value = -618.2665999999999
$ws.Range("N102").Value = -755541.4
$ws.Range("H122").Value = 1700
$ws.Range("I122").Value = 1050
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 3150
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -700
$ws.Range("N122").Value = -13900
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H132").Value = 2774.9644
$ws.Range("I132").Value = 2245.05
$ws.Range("J132").Value = 4099.75
$ws.Range("K132").Value = 6735.150000000001
$ws.Range("L132").Value = 12299.25
$ws.Range("M132").Value = -4205.150000000001
$ws.Range("N132").Value = -17359.25
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1522.921
$ws.Range("I132").Value = 1264.1482
$ws.Range("J132").Value = 2158.0908
$ws.Range("K132").Value = 3792.4446
$ws.Range("L132").Value = 6474.2724
$ws.Range("M132").Value = -1262.4446
$ws.Range("N132").Value = -11534.2724
$ws.Range("H136").Value = 1897

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1177.2565
$ws.Range("I132").Value = 777.5806
$ws.Range("J132").Value = 2726
$ws.Range("K132").Value = 2332.7418
$ws.Range("L132").Value = 8178
$ws.Range("M132").Value = 197.2582000000002
$ws.Range("N132").Value = -13238
$ws.Range("H136").Value = 1014.2857
$ws.Range("J136").Value = 2500
$ws.Range("L136").Value = 7500
$ws.Range("N136").Value = -12600
